$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) holds text-formatted numbers (e.g. "28.455.38", "86.70").
# Force text format before writing so Excel does not coerce them to numbers
# and silently drop formatting such as trailing zeros.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '28.455.38'
$ws.Range('E2').Value = '  +1.18%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.864.64'
$ws.Range('E3').Value = '  +1.48%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '325.01'
$ws.Range('E5').Value = '  -0.21%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.006'
$ws.Range('E6').Value = '  +0.02%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4556'
$ws.Range('E7').Value = '  -1.76%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3839'
$ws.Range('E8').Value = '  -0.51%  '
$ws.Range('E9').Value = '  -0.13%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.9888'
$ws.Range('E10').Value = '  +2.80%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '21.55'
$ws.Range('E11').Value = '  -2.28%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.879.22'
$ws.Range('E12').Value = '  +7.49%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.915'
$ws.Range('E13').Value = '  +0.77%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.644'
$ws.Range('E14').Value = '  -0.75%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.06957'
$ws.Range('E15').Value = '  +0.93%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '86.70'
$ws.Range('E16').Value = '  -2.04%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.007'
$ws.Range('E17').Value = '  +0.11%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000009950'
$ws.Range('E18').Value = '  +0.08%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '16.65'
$ws.Range('E19').Value = '  -0.17%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.007'
$ws.Range('E20').Value = '  +0.25%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '28.483.95'
$ws.Range('E21').Value = '  +1.14%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.254'
$ws.Range('E22').Value = '  -0.77%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.91'
$ws.Range('E23').Value = '  -0.97%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.099'
$ws.Range('E24').Value = '  +0.08%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.094.03'
$ws.Range('E25').Value = '  +4.56%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '153.43'
$ws.Range('E26').Value = '  -0.62%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '19.14'
$ws.Range('E27').Value = '  +0.05%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '5.663'
$ws.Range('E28').Value = '  -1.19%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.940'
$ws.Range('E29').Value = '  -1.37%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '117.56'
$ws.Range('E30').Value = '  -0.90%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.09277'
$ws.Range('E31').Value = '  +0.21%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.9108'
$ws.Range('E32').Value = '  -2.08%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.268'
$ws.Range('E33').Value = '  -0.22%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.320'
$ws.Range('E34').Value = '  -0.36%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.294'
$ws.Range('E35').Value = '  -1.14%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.05722'
$ws.Range('E36').Value = '  -1.50%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.141'
$ws.Range('E37').Value = '  -0.23%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02057'
$ws.Range('E38').Value = '  -2.67%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '7.677'
$ws.Range('E39').Value = '  -1.06%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.5568'
$ws.Range('E40').Value = '  -0.23%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.1771'
$ws.Range('E41').Value = '  +0.57%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '9.646'
$ws.Range('E42').Value = '  -2.07%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.07098'
$ws.Range('E43').Value = '  -1.95%  '
$ws.Range('E44').Value = '  -0.47%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.5241'
$ws.Range('E45').Value = '  -0.38%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.145'
$ws.Range('E46').Value = '  +0.99%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.122'
$ws.Range('E47').Value = '  -0.90%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.814'
$ws.Range('E48').Value = '  -0.96%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '111.73'
$ws.Range('E49').Value = '  -1.81%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.411'
$ws.Range('E50').Value = '  +3.66%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.005'
$ws.Range('E51').Value = '  -0.02%  '
